$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H3 to a new value ("Standalone, WC, test") which will be added
# as a new shared string.
$ws.Range("H3").Value = "Standalone, WC, test"

# Update the selection/active cell shown in the sheet view to H6.
$ws.Range("H6").Select()
